$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 19

# Capture existing ExpPoints values (column C, rows 2-19) before they get
# overwritten, so they can be moved over to the new column G.
$expPoints = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $expPoints[$r] = $ws.Cells.Item($r, 3).Value2
}

# Set the new header cells' text. The previous "ExpPoints" header (C1) now
# becomes "WIN"; four new headers are inserted after it (TOP2, TOP4,
# RELEGATION), and the ExpPoints header moves out to the new last column G.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"

# Apply the same bold/centered/bordered header formatting (taken from the
# existing A1 header cell) to the newly added D1:G1 header cells.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# PasteSpecial(xlPasteFormats) only copies formatting, but make sure the
# text values are still correct afterwards.
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"

# For each data row: clear out the new placeholder columns C,D,E,F, and move
# the previously-read ExpPoints value into the new column G.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 7).Value = $expPoints[$r]
}
